$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) with the same header style as the other
# header cells (e.g. G1: bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new column's data rows (plain numeric, no special style).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
